# Updated the test cases to include write to excel utility
#
# Inserts two new worksheets - "MTTS" and "Territory" - between the
# existing "Dispatch_Process" and "RS_6299_RS_6300" sheets, and
# populates them with their rule-name data. Adding "Territory" last
# makes it the active/selected sheet, matching the target workbook.

$wb = $excel.ActiveWorkbook

# --- MTTS sheet: inserted right after Dispatch_Process ---
$dispatch = $wb.Worksheets.Item("Dispatch_Process")
$mtts = $wb.Worksheets.Add($null, $dispatch)
$mtts.Name = "MTTS"

$mtts.Range("A1").Value = "MTTS Name"
# Written oldest-first so the shared-string table / cell order matches
# the source data (rows display newest-first, top to bottom).
$mtts.Range("A4").Value = "MTTSRule_Wed Oct 04 2017 14:13:01 GMT+0530 (IST)"
$mtts.Range("A3").Value = "MTTSRule_Wed Oct 04 2017 14:17:09 GMT+0530 (IST)"
$mtts.Range("A2").Value = "MTTSRule_Wed Oct 04 2017 14:24:25 GMT+0530 (IST)"

# --- Territory sheet: inserted right after MTTS ---
$territory = $wb.Worksheets.Add($null, $mtts)
$territory.Name = "Territory"

$territory.Range("A1").Value = "Territory Name"
$territory.Range("A2").Value = "TerritoryRule_Wed Oct 04 2017 14:41:41 GMT+0530 (IST)"
